# Swap the match-record data (columns B:AC) between paired rows.
# The "id" column (A) is a plain sequential row index and is left untouched;
# everything else (match id, teams, scores, odds, ...) is exchanged between
# the two rows of each pair, matching the source-data re-ordering in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(7, 8),
    @(31, 32),
    @(43, 44),
    @(61, 62),
    @(77, 78),
    @(90, 91),
    @(214, 215)
)

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $range1 = $ws.Range("B" + $row1 + ":AC" + $row1)
    $range2 = $ws.Range("B" + $row2 + ":AC" + $row2)

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}
